# Add: CherryMX Switch and Rectifier Diode Components.
# This script reproduces, via Excel COM-interop calls, the addition of two
# new component sections ("Header Pins (16)" and "Rectifier Diode") to the
# right-hand block (columns L:S, plus X) of the Components sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Components")

# ---------------------------------------------------------------------
# 1. Copy the formatting of an existing analogous block down onto the
#    rows/columns we are about to populate, so the new cells pick up the
#    same styles (header banner / column-title / data-row styles) already
#    used elsewhere in this L:S side table.
# ---------------------------------------------------------------------

# "Header Pins (16)" section -> rows 23,24,25 (pattern copied from rows 19,20,21)
$ws.Range("L19:S21").Copy()
$ws.Range("L23").PasteSpecial(-4122)

# "Rectifier Diode" section -> rows 26,27,28 (pattern copied from rows 19,20,21 again)
$ws.Range("L19:S21").Copy()
$ws.Range("L26").PasteSpecial(-4122)

$ws.Range("Calculation").Select()
$excel.CutCopyMode = 0

# Single extra hyperlink cell X28 picks up the same "plain link" style used
# for similar single hyperlink cells elsewhere (e.g. M13 / M14 / I32).
$ws.Range("M13").Copy()
$ws.Range("X28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Populate the cell values. The order below matters: it controls the
#    order new entries are appended to the shared-string table.
# ---------------------------------------------------------------------

$ws.Range("L23").Value = "Header Pins (16)"

$ws.Range("M25").Value = "https://www.digikey.com.au/product-detail/en/amphenol-icc-fci/67996-416HLF/609-3220-ND/1878538"
$ws.Range("L25").Value = "67996-416HLF"

$ws.Range("M28").Value = "https://www.digikey.com.au/product-detail/en/vishay-general-semiconductor-diodes-division/ES07D-GS08/ES07D-GS08CT-ND/3104461"
$ws.Range("L26").Value = "Rectifier Diode"
$ws.Range("L28").Value = "ES07D-GS08"
$ws.Range("X28").Value = "https://www.digikey.com.au/product-detail/en/vishay-general-semiconductor-diodes-division/ES07D-GS08/ES07D-GS08CT-ND/3104461"

# Column headers (Name / Supplier Link / Datasheet / ORDERED) for each new section
$ws.Range("L24").Value = "Name"
$ws.Range("M24").Value = "Supplier Link"
$ws.Range("P24").Value = "Datasheet"
$ws.Range("S24").Value = "ORDERED"

$ws.Range("L27").Value = "Name"
$ws.Range("M27").Value = "Supplier Link"
$ws.Range("P27").Value = "Datasheet"
$ws.Range("S27").Value = "ORDERED"

# ---------------------------------------------------------------------
# 3. Merge the header / link / datasheet cells, matching the layout used
#    by the other sections of this table.
# ---------------------------------------------------------------------

$ws.Range("L23:S23").Merge()
$ws.Range("M24:O24").Merge()
$ws.Range("P24:R24").Merge()
$ws.Range("M25:O25").Merge()
$ws.Range("P25:R25").Merge()

$ws.Range("L26:S26").Merge()
$ws.Range("M27:O27").Merge()
$ws.Range("P27:R27").Merge()
$ws.Range("M28:O28").Merge()
$ws.Range("P28:R28").Merge()

# ---------------------------------------------------------------------
# 4. Hyperlinks for the supplier-link cells.
# ---------------------------------------------------------------------

$ws.Hyperlinks.Add($ws.Range("M25"), "https://www.digikey.com.au/product-detail/en/amphenol-icc-fci/67996-416HLF/609-3220-ND/1878538")
$ws.Hyperlinks.Add($ws.Range("M28"), "https://www.digikey.com.au/product-detail/en/vishay-general-semiconductor-diodes-division/ES07D-GS08/ES07D-GS08CT-ND/3104461")
$ws.Hyperlinks.Add($ws.Range("X28"), "https://www.digikey.com.au/product-detail/en/vishay-general-semiconductor-diodes-division/ES07D-GS08/ES07D-GS08CT-ND/3104461")

# ---------------------------------------------------------------------
# 5. View / selection state: the Components sheet becomes the active tab,
#    scrolled down a bit with the last worked-on cell selected; the
#    "Atmega32U2 Pinout" sheet loses the active-tab flag it used to have.
# ---------------------------------------------------------------------

$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("L29").Select()
